# Logboek - week 7: add Tuesday (Dinsdag) entries, push Thursday (Donderdag)
# entry down a row, and clear the now-unused numbering in rows 11-14.
#
# Commit message: "Mogelijk gemaakt om bestanden weg te laten schijven"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 7")

# Row 7: day label changes from Donderdag to Dinsdag (date 41681 = 2014-02-11,
# a Tuesday, already in the sheet - only the label was wrong/changed)
$ws.Range("A7").Value = "Dinsdag"

# Row 9: new log entry for Dinsdag - 10:45 to 12:00
$ws.Range("C9").Value = 0.44791666666666669
$ws.Range("D9").Value = 0.5
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = "Mogelijk gemaakt om bestanden weg te laten schijven"
$ws.Rows.Item(9).RowHeight = 30

# Row 10: new log entry for Donderdag (2014-02-13) - 10:30 to 11:00
$ws.Range("A10").Value = "Donderdag"
$ws.Range("B10").Value = 41683
$ws.Range("C10").Value = 0.4375
$ws.Range("D10").Value = 0.45833333333333331
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = "Jquery gedownload"

# Rows 11-14: the Id numbering that used to live here moved up into rows 9-10,
# so clear it out
$ws.Range("E11").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("E14").ClearContents()

# Update the view state to match where the user ended up (selection on F11)
$ws.Activate()
$ws.Range("F11").Select()
